$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal value that Excel would otherwise auto-convert
# (e.g. "211.76") as plain text, matching the original inline-string cells
# with no number-format / style changes.
function Set-TextLiteral($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $escaped = $val -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

$ws.Range("D2").Value = '28.402.47'
$ws.Range("E2").Value = '  -0.30%  '
$ws.Range("D3").Value = '1.571.10'
$ws.Range("E3").Value = '  -0.02%  '
$ws.Range("E4").Value = '  -0.04%  '
Set-TextLiteral 5 4 '211.76'
$ws.Range("E5").Value = '  -0.06%  '
Set-TextLiteral 6 4 '0.492'
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("E7").Value = '  +0.01%  '
Set-TextLiteral 8 4 '44.40'
$ws.Range("E8").Value = '  -3.73%  '
Set-TextLiteral 9 4 '23.65'
$ws.Range("E9").Value = '  -1.84%  '
$ws.Range("E10").Value = '  -0.56%  '
$ws.Range("E11").Value = '  -0.57%  '
Set-TextLiteral 12 4 '0.0893'
$ws.Range("E12").Value = '  +1.28%  '
$ws.Range("D13").Value = '1.794.96'
$ws.Range("E13").Value = '  -0.05%  '
$ws.Range("D14").Value = '1.576.13'
$ws.Range("E14").Value = '  +0.32%  '
$ws.Range("E15").Value = '  -0.07%  '
$ws.Range("D16").Value = '28.415.55'
$ws.Range("E16").Value = '  -0.14%  '
$ws.Range("E17").Value = '  -1.25%  '
$ws.Range("E18").Value = '  -0.80%  '
Set-TextLiteral 19 4 '228.29'
$ws.Range("E19").Value = '  +0.48%  '
$ws.Range("E20").Value = '  +0.43%  '
$ws.Range("E21").Value = '  -1.11%  '
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("E23").Value = '  +1.38%  '
Set-TextLiteral 24 4 '8.98'
$ws.Range("E24").Value = '  -1.53%  '
Set-TextLiteral 25 4 '2.06'
$ws.Range("E25").Value = '  +0.81%  '
Set-TextLiteral 26 4 '151.29'
$ws.Range("E26").Value = '  +0.15%  '
$ws.Range("E27").Value = '  -0.53%  '
$ws.Range("E28").Value = '  -0.44%  '
Set-TextLiteral 29 4 '6.35'
$ws.Range("E29").Value = '  -1.32%  '
$ws.Range("E30").Value = '  -0.03%  '
Set-TextLiteral 31 4 '0.0481'
$ws.Range("E31").Value = '  +3.67%  '
$ws.Range("E32").Value = '  -2.11%  '
$ws.Range("E33").Value = '  -0.80%  '
$ws.Range("E34").Value = '  -1.44%  '
$ws.Range("D35").Value = '1.379.70'
$ws.Range("E35").Value = '  -0.77%  '
$ws.Range("E36").Value = '  +4.34%  '
$ws.Range("E37").Value = '  -2.14%  '
$ws.Range("E38").Value = '  +0.27%  '
$ws.Range("E39").Value = '  +0.91%  '
$ws.Range("E40").Value = '  -1.52%  '
$ws.Range("E41").Value = '  -1.98%  '
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("E43").Value = '  +2.51%  '
Set-TextLiteral 44 4 '0.787'
$ws.Range("E44").Value = '  -0.68%  '
Set-TextLiteral 45 4 '0.0473'
$ws.Range("E45").Value = '  +0.68%  '
$ws.Range("E46").Value = '  -4.25%  '
Set-TextLiteral 47 4 '62.24'
$ws.Range("E47").Value = '  -1.17%  '
$ws.Range("E48").Value = '  -6.16%  '
$ws.Range("D49").Value = '1.707.89'
$ws.Range("E49").Value = '  +0.01%  '
Set-TextLiteral 50 4 '2.14'
$ws.Range("E50").Value = '  +0.03%  '
Set-TextLiteral 51 4 '85.26'
$ws.Range("E51").Value = '  -0.72%  '
